# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada
# das linhas com total das tabelas.
#
# For each of the first five "year" tables the E1 header cell holds a
# leftover numeric value (688.0824318511211) instead of the text label
# that belongs there ("2050", or "2041-2050" on the incremental-power
# sheet). We fix the label and then drop the trailing "Total" row that
# shouldn't be part of these tables.

$wb = $excel.ActiveWorkbook

function Set-TextLabel {
    # Writes $text into $range as a genuine text value (not a number),
    # while keeping the cell's existing style (font/border/alignment)
    # intact. A direct `$range.Value = "2050"` would be auto-coerced to
    # a number by Excel because the string looks numeric, and prefixing
    # it with a quote directly on the target cell would tattoo a new
    # "quote prefixed" style onto it. Routing the text through a scratch
    # cell and pasting values-only, then formats-only from a cell that
    # already carries the style we want, avoids both problems.
    param(
        $Worksheet,
        [string]$Address,
        [string]$Text,
        [string]$StyleSourceAddress
    )
    # NB: called positionally everywhere below -- this runtime's PS
    # parser does not bind `-Name value` style named arguments.

    $target = $Worksheet.Range($Address)
    $styleSource = $Worksheet.Range($StyleSourceAddress)
    $scratch = $Worksheet.Range("ZZ1")

    $scratch.Value = "'" + $Text
    $scratch.Copy()
    $target.PasteSpecial(-4163)   # xlPasteValues

    $styleSource.Copy()
    $target.PasteSpecial(-4122)   # xlPasteFormats

    $scratch.Clear()
    $Worksheet.Application.CutCopyMode = $false
}

# --- Sheet 1: "Potencia Acumulada - SIN (MW)" ------------------------
$ws = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

# --- Sheet 2: "Geracao Periodo Medio (MWMed)" -------------------------
$ws = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

# --- Sheet 3: "Atendimento a Ponta(MW)" -------------------------------
$ws = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

# --- Sheet 4: "Potencia Incremental - SIN(MW)" ------------------------
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws "E1" "2041-2050" "D1"
$ws.Rows.Item(13).Delete()

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" -----------------------------
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws "E1" "2050" "D1"
# this table never had a Total row, nothing to delete

# --- Sheet 6: "Custo Total (bilhoes de R$)" ---------------------------
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Rows.Item(4).Delete()
